# Update "想去人数" (F column) counts for several rows across the four sheets
# to match the latest scraped data (commit: Update gh-pages to output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F3"  = 419
        "F6"  = 3814
        "F8"  = 2537
        "F9"  = 72
        "F10" = 3086
        "F12" = 530
        "F13" = 2294
        "F17" = 445
        "F18" = 2
        "F20" = 203
        "F23" = 356
        "F24" = 649
        "F25" = 1401
        "F30" = 3
        "F32" = 42
        "F33" = 4250
        "F34" = 3931
        "F35" = 73
        "F36" = 4
        "F38" = 1117
        "F39" = 2
        "F43" = 159
        "F45" = 95
        "F48" = 55
    }
    "演出" = @{
        "F15" = 204
    }
    "本地生活" = @{
        "F2" = 1027
        "F3" = 144
        "F4" = 2271
    }
    "全部类型" = @{
        "F3"  = 1027
        "F4"  = 144
        "F5"  = 419
        "F10" = 3814
        "F12" = 2537
        "F13" = 72
        "F14" = 3086
        "F15" = 530
        "F16" = 2294
        "F20" = 445
        "F24" = 356
        "F25" = 649
        "F26" = 1401
        "F32" = 42
        "F35" = 4251
        "F36" = 3931
        "F37" = 73
        "F38" = 1117
        "F45" = 159
        "F46" = 95
        "F48" = 55
        "F49" = 204
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
